# Updates cryptocurrency price/volume data in the worksheet to match
# the latest scrape. Numeric-looking "price" strings (e.g. "0.999",
# "233.01") must stay as TEXT, matching the source feed which renders
# prices/volumes as formatted strings rather than numbers -- the sheet
# already stores values like "1.00" and "621.06" as text for the same
# reason. We briefly mark such cells as Text (NumberFormat "@") before
# assigning the value so Excel does not auto-convert them to numbers,
# then restore the "Normal" style so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '90.508.08'
$ws.Range("E2").Value = '  -0.38%  '
$ws.Range("D3").Value = '3.100.87'
$ws.Range("E3").Value = '  -1.90%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.32%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.01'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +8.17%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '621.88'
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = '  -4.17%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.369'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.07%  '
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("D10").Value = '3.098.09'
$ws.Range("E10").Value = '  -1.89%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.739'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.03%  '
$ws.Range("E12").Value = '  -2.80%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000252'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.51%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.98'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.75%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.50'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.55%  '
$ws.Range("D16").Value = '90.194.68'
$ws.Range("E16").Value = '  -0.38%  '
$ws.Range("D17").Value = '3.661.50'
$ws.Range("E17").Value = '  -1.87%  '
$ws.Range("D18").Value = '3.095.00'
$ws.Range("E18").Value = '  -1.47%  '
$ws.Range("E19").Value = '  +3.87%  '
$ws.Range("E20").Value = '  +5.35%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.07'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.71%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '437.28'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -6.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.57'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.94'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.66%  '
$ws.Range("B25").Value = 'LEO'
$ws.Range("C25").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '7.56'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.52%  '
$ws.Range("B26").Value = 'NEARProtocol'
$ws.Range("C26").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '5.69'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.35%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '89.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.03%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.10'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.37%  '
$ws.Range("D29").Value = '3.246.70'
$ws.Range("E29").Value = '  -2.07%  '
$ws.Range("E30").Value = '  -0.12%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.33'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.50%  '
$ws.Range("E32").Value = '  -0.64%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.73%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.197'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +9.73%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '25.87'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.97%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.153'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +7.89%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.82'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.26'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.20%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '503.37'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.25%  '
$ws.Range("E40").Value = '  -1.07%  '
$ws.Range("E41").Value = '  -2.14%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0889'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.04%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '22.18'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.14%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.406'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.26%  '
$ws.Range("E45").Value = '  +0.01%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.45'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +54.67%  '
$ws.Range("E47").Value = '  -3.55%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.690'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.89%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '152.56'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.31%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '44.90'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.35%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.34'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.58%  '
